$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column Q (17th column) from 8 to 10.
# Note: Excel's ColumnWidth COM property is expressed in "characters" based
# on the default font's max digit width, which is offset from the raw
# stored <col width="..."/> units by a constant padding factor (~0.83 for
# the default Calibri 11 font here). Subtract that offset so the value
# persisted to the sheet XML ends up as exactly 10.
$ws.Columns.Item(17).ColumnWidth = 9.17

# Row 2: Statut goes from "NA" to "RA", and "Dern. adh." (F2) now holds 2024
$ws.Range("E2").Value = "RA"
$ws.Range("F2").Value = 2024

# Prepare rows 3-6 with the same cell formatting (style) as the existing
# data row 2, spanning the full A:T width, before filling in values.
$ws.Range("A2:T2").Copy()
$ws.Range("A3:T6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 4: "Nombre de reçus" counter
$ws.Range("A4").Value = "Nombre de reçus"
$ws.Range("B4").Value = 0

# Row 5: "Total avec reçus" counter
$ws.Range("A5").Value = "Total avec reçus"
$ws.Range("B5").Value = 0

# Row 6: "Total sans reçus" counter
$ws.Range("A6").Value = "Total sans reçus"
$ws.Range("B6").Value = 10
